# modification v 3.3 suivi
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate()

# Mark the six tasks in rows 5-10 as "A faire" (to-do) by copying their
# label from column B into column C, keeping the same text (reuses the
# existing shared string) and picking up the yellow "A faire" highlight
# formatting already used elsewhere in the sheet (e.g. C39/C40) via a
# format-only paste.
for ($r = 5; $r -le 10; $r++) {
    $ws.Range("C39").Copy() | Out-Null
    $ws.Range("C$r").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
    $ws.Range("C$r").Value2 = $ws.Range("B$r").Value2
}
$excel.CutCopyMode = $false

# Scroll the view back to the top of the sheet and move the active
# selection to C10, matching where the edit left off.
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("C10").Select()
